$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G9").Value = "Acierto"
$ws.Range("H9").Value = 2.25
